$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H105").Value = 12000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 12000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 12000
$ws.Range("N105").Value = -18988

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

$ws.Range("H74").Value = 2860.2
$ws.Range("I74").Value = 2622.4443
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 2622.4443
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -1748.4443
$ws.Range("N74").Value = -6748

$ws.Range("H77").Value = 2860.2
$ws.Range("I77").Value = 2622.4443
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 13112.2215
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -8744.2215
$ws.Range("N77").Value = -33736

$ws.Range("H95").Value = 23201.285
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 23201.285
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 23201.285
$ws.Range("N95").Value = -28693.285

$ws.Range("H104").Value = 14306
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 14306
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 14306
$ws.Range("N104").Value = -21294

$ws.Range("H105").Value = 42684.5
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 42684.5
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 42684.5
$ws.Range("N105").Value = -49672.5

$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()

$ws.Range("H122").Value = 1703.5
$ws.Range("I122").Value = 450
$ws.Range("J122").Value = 2957
$ws.Range("K122").Value = 1350
$ws.Range("L122").Value = 8871
$ws.Range("M122").Value = 1100
$ws.Range("N122").Value = -13771

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()

$ws.Range("H5").Value = 377
$ws.Range("I5").Value = 206.66667
$ws.Range("J5").Value = 581.4
$ws.Range("K5").Value = 206.66667
$ws.Range("L5").Value = 581.4
$ws.Range("M5").Value = -93.66667000000001
$ws.Range("N5").Value = -807.4

$ws.Range("H80").Value = 830.6667
$ws.Range("I80").Value = 676.6
$ws.Range("J80").Value = 1023.25
$ws.Range("K80").Value = 676.6
$ws.Range("L80").Value = 1023.25
$ws.Range("M80").Value = 321.4
$ws.Range("N80").Value = -3019.25

$ws.Range("H83").Value = 830.6667
$ws.Range("I83").Value = 676.6
$ws.Range("J83").Value = 1023.25
$ws.Range("K83").Value = 3383
$ws.Range("L83").Value = 5116.25
$ws.Range("M83").Value = 1609
$ws.Range("N83").Value = -15100.25

$ws.Range("H86").Value = 450
$ws.Range("I86").Value = 450
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 450
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 673

$ws.Range("H89").Value = 450
$ws.Range("I89").Value = 450
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 2250
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 3366

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H94").Value = 7933.3335
$ws.Range("I94").Value = 8250
$ws.Range("J94").Value = 7300
$ws.Range("K94").Value = 8250
$ws.Range("L94").Value = 7300
$ws.Range("M94").Value = -7799
$ws.Range("N94").Value = -8202

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H104").Value = 59950
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 59950
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 59950
$ws.Range("N104").Value = -66938

$ws.Range("H122").Value = 6664
$ws.Range("I122").Value = 6664
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 19992
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -17542

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 1000000
$ws.Range("I4").Value = 1000000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1000000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -999887

$ws.Range("H7").Value = 3497.25
$ws.Range("I7").Value = 3139.7144
$ws.Range("J7").Value = 6000
$ws.Range("K7").Value = 3139.7144
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = -3027.7144
$ws.Range("N7").Value = -6224

$ws.Range("H28").Value = 1000000
$ws.Range("I28").Value = 1000000
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 1000000
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -999768

$ws.Range("H37").Value = 1000000
$ws.Range("I37").Value = 1000000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 1000000
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -999893

$ws.Range("H40").Value = 2852.3333
$ws.Range("I40").Value = 2852.3333
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2852.3333
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2716.3333

$ws.Range("H61").Value = 1059.6
$ws.Range("I61").Value = 1098.25
$ws.Range("J61").Value = 905
$ws.Range("K61").Value = 1098.25
$ws.Range("L61").Value = 905
$ws.Range("M61").Value = -896.25
$ws.Range("N61").Value = -1309

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H105").Value = 26000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 26000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 26000
$ws.Range("N105").Value = -32988

$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H113").Value = 1059.6
$ws.Range("I113").Value = 1098.25
$ws.Range("J113").Value = 905
$ws.Range("K113").Value = 1098.25
$ws.Range("L113").Value = 905
$ws.Range("M113").Value = 1071.75
$ws.Range("N113").Value = -5245

$ws.Range("H126").Value = 3497.25
$ws.Range("I126").Value = 3139.7144
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 9419.143199999999
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -6949.143199999999
$ws.Range("N126").Value = -22940

$ws.Range("H139").Value = 45197.25
$ws.Range("I139").Value = 45789
$ws.Range("J139").Value = 45000
$ws.Range("K139").Value = 45789
$ws.Range("L139").Value = 45000
$ws.Range("M139").Value = -40649
$ws.Range("N139").Value = -55280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 18500
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 18500
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 18500
$ws.Range("N103").Value = -20844

$ws.Range("H104").Value = 10308.667
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 10308.667
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 10308.667
$ws.Range("N104").Value = -17296.667

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H110").Value = 25000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 25000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 25000
$ws.Range("N110").Value = -33180

$ws.Range("H135").Value = 61542
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 61542
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 61542
$ws.Range("N135").Value = -71682

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
